$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their original text formatting
# so that numeric-looking strings (e.g. "5.00", "1.30") are not silently
# converted into numbers and stripped of trailing zeros by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '72.013.14'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '2.683.19'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '598.12'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('D6').Value = '174.45'
$ws.Range('E6').Value = '  -3.69%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = '2.681.93'
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('D10').Value = '0.167'
$ws.Range('E10').Value = '  -5.86%  '
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('D12').Value = '0.356'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').Value = '5.00'
$ws.Range('E13').Value = '  -1.56%  '
$ws.Range('D14').Value = '3.173.72'
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0000184'
$ws.Range('E15').Value = '  -5.53%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '71.913.80'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '26.17'
$ws.Range('D18').Value = '2.685.21'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').Value = '12.16'
$ws.Range('E19').Value = '  +4.48%  '
$ws.Range('D20').Value = '8.13'
$ws.Range('E20').Value = '  +2.90%  '
$ws.Range('D21').Value = '372.28'
$ws.Range('E21').Value = '  -3.54%  '
$ws.Range('D22').Value = '4.17'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = '2.01'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '72.37'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = '4.34'
$ws.Range('E26').Value = '  -2.29%  '
$ws.Range('D27').Value = '9.80'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').Value = '2.820.93'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '0.0₃0976'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '8.04'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').Value = '502.94'
$ws.Range('E32').Value = '  -8.54%  '
$ws.Range('D33').Value = '1.30'
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '164.06'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('D37').Value = '19.63'
$ws.Range('E37').Value = '  +1.18%  '
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '1.38'
$ws.Range('E39').Value = '  -2.92%  '
$ws.Range('D40').Value = '0.109'
$ws.Range('E40').Value = '  -4.83%  '
$ws.Range('E41').Value = '  -5.42%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').Value = '5.01'
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '2.57'
$ws.Range('E44').Value = '  -2.50%  '
$ws.Range('B45').Value = 'PolygonEcosystemToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D45').Value = '0.332'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('D46').Value = '156.46'
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').Value = '39.42'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').Value = '0.563'
$ws.Range('E48').Value = '  +3.83%  '
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('E51').Value = '  -0.47%  '
